# Update the "取得日時" (retrieved-at) timestamp column on the "ランサーズ" sheet.
# All existing data rows (2-10) get their timestamp refreshed to the latest
# scrape time, as recorded in the commit message: 2025-12-12 12:51 JST.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-12 12:51:58"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 2
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and "$($cell.Value)" -ne "") {
        $cell.Value = $newTimestamp
    }
}
